$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) columns for rows 2-4 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 100.65971175
$wsSchedule.Range("F2").Value = 1.972559509112287
$wsSchedule.Range("E3").Value = 451.2848925
$wsSchedule.Range("F3").Value = 21.70682503607504
$wsSchedule.Range("E4").Value = -2.813537999999978
$wsSchedule.Range("F4").Value = -0.08270246913580183

# --- Sheet "Detailed": update Price column (B) and Type column (C) for rows 29-94 ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B29").Value = -7
$wsDetailed.Range("B30").Value = -12.20552
$wsDetailed.Range("B31").Value = -7.46266
$wsDetailed.Range("B32").Value = -6.8
$wsDetailed.Range("B33").Value = -6.32
$wsDetailed.Range("B34").Value = -5.41
$wsDetailed.Range("B35").Value = -4.77891
$wsDetailed.Range("B36").Value = 5.16051
$wsDetailed.Range("B37").Value = 6.2382
$wsDetailed.Range("B38").Value = 6.99021
$wsDetailed.Range("B39").Value = 18.34255
$wsDetailed.Range("B40").Value = 36.2
$wsDetailed.Range("B41").Value = 45.74218
$wsDetailed.Range("B42").Value = 49.42693
$wsDetailed.Range("B44").Value = 46.19026
$wsDetailed.Range("B46").Value = 47.91617
$wsDetailed.Range("B52").Value = 56.98
$wsDetailed.Range("B53").Value = 56.78
$wsDetailed.Range("B55").Value = 48.75983
$wsDetailed.Range("B56").Value = 49.44748
$wsDetailed.Range("B57").Value = 50.33883
$wsDetailed.Range("B59").Value = 57.76254
$wsDetailed.Range("B60").Value = 57.06003
$wsDetailed.Range("B61").Value = 58.81169
$wsDetailed.Range("B62").Value = 58.02194
$wsDetailed.Range("B64").Value = 36.06
$wsDetailed.Range("B65").Value = 36.06
$wsDetailed.Range("B66").Value = 5.16569
$wsDetailed.Range("B67").Value = 0.64744
$wsDetailed.Range("B69").Value = -5.50985
$wsDetailed.Range("B70").Value = -6.14712
$wsDetailed.Range("B71").Value = -7.02049
$wsDetailed.Range("B72").Value = -7.15517
$wsDetailed.Range("B73").Value = -7.67966
$wsDetailed.Range("B74").Value = -7.68692
$wsDetailed.Range("B75").Value = -7.23621
$wsDetailed.Range("B76").Value = -6.19058
$wsDetailed.Range("B77").Value = -6.35843
$wsDetailed.Range("B78").Value = -6.06857
$wsDetailed.Range("B79").Value = -5.72243
$wsDetailed.Range("B80").Value = -5.51
$wsDetailed.Range("B81").Value = -2.53338
$wsDetailed.Range("B82").Value = 0
$wsDetailed.Range("B83").Value = -2.54896
$wsDetailed.Range("B84").Value = -2.40408
$wsDetailed.Range("B85").Value = 4.79045
$wsDetailed.Range("B86").Value = 20.78901
$wsDetailed.Range("B87").Value = 54.66255
$wsDetailed.Range("B88").Value = 57.85505
$wsDetailed.Range("B89").Value = 71.91437000000001
$wsDetailed.Range("B90").Value = 68.63066999999999
$wsDetailed.Range("B91").Value = 66.50706
$wsDetailed.Range("B92").Value = 63.83273
$wsDetailed.Range("B93").Value = 61.9517
$wsDetailed.Range("B94").Value = 58.95387

$wsDetailed.Range("C32").Value = "historical"
$wsDetailed.Range("C33").Value = "historical"
